$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: recalc similarity_score for "formats for loading a text file" (deepseek1.5)
$ws.Range("F2").Value = -0.02518575824797153

# Row 3: fill in Expected Answer + recalc similarity_score for "tracks in one ODF" (llama3.2:latest)
$ws.Range("E3").Value = "200 tracks can be defined in one ODF."
$ws.Range("F3").Value = 0.7644559144973755

# Row 4: "curves can I load in one go" (deepseek1.5) - Expected Answer text unchanged, score unchanged
$ws.Range("E4").Value = "450 curves can be loaded in one go."
$ws.Range("F4").Value = 0.9253911972045898

# Row 5: "curves can I load in one go" (llama3.2:latest) - Expected Answer text unchanged, score unchanged
$ws.Range("E5").Value = "450 curves can be loaded in one go."
$ws.Range("F5").Value = 0.9253911972045898

# Row 6: fill in Expected Answer + recalc similarity_score for "max number of characters" (llama3.2:latest)
$ws.Range("E6").Value = "250 / 32000 (varies per text type) are the maximum number of characters in a single text entry."
$ws.Range("F6").Value = 0.913228452205658
